$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Duplicate the "Add Panels" sheet (in its ORIGINAL, unedited state)
#    to the end of the workbook and rename it "Sheet2". This preserves
#    the original "FIRECLASS 64-2" label (row 8) on the copy.
# ------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("Add Panels")
$srcSheet.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Sheet2"

# Add a new data row (row 9) to the duplicated sheet for the "FC708D"
# panel -- values first, then copy row 8's formatting onto row 9 so the
# styles (incl. the odd quote-prefixed numeric styles) come out right.
$newSheet.Cells.Item(9, 1).Value = "FC708D"
$newSheet.Cells.Item(9, 2).Value = "Node1"
$newSheet.Cells.Item(9, 3).Value = "CPU 800"
$newSheet.Cells.Item(9, 4).Value = "PFI"
$newSheet.Cells.Item(9, 5).Value = 14
$newSheet.Cells.Item(9, 6).Value = 0.276
$newSheet.Cells.Item(9, 7).Value = 0.426
$newSheet.Cells.Item(9, 9).Value = "MPM800"
$newSheet.Cells.Item(9, 10).Value = "Miscellaneous"
$newSheet.Cells.Item(9, 11).Value = "MPM800-1"
$newSheet.Cells.Item(9, 12).Value = 0.022
$newSheet.Cells.Item(9, 13).Value = 0.03
$newSheet.Cells.Item(9, 14).Value = 0.298
$newSheet.Cells.Item(9, 15).Value = 0.456
$newSheet.Cells.Item(9, 16).Value = "5V"
$newSheet.Cells.Item(9, 17).Value = 0.286
$newSheet.Cells.Item(9, 18).Value = 0.445

$newSheet.Range("A8:G8").Copy()
$newSheet.Range("A9:G9").PasteSpecial(-4122)
$newSheet.Range("I8:R8").Copy()
$newSheet.Range("I9:R9").PasteSpecial(-4122)

$newSheet.Cells.Select()

# ------------------------------------------------------------------
# 2) Edit the original "Add Panels" sheet: row 8's panel name changes
#    from "FIRECLASS 64-2" to "FC64-2", and a brand-new row 9 is
#    appended for the "FC708D" panel (its CPU-Type cell is blank).
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Add Panels")

$ws.Cells.Item(8, 1).Value = "FC64-2"

$ws.Cells.Item(9, 1).Value = "FC708D"
$ws.Cells.Item(9, 2).Value = "Node1"
$ws.Cells.Item(9, 3).Formula = "'"
$ws.Cells.Item(9, 4).Value = "PFI"
$ws.Cells.Item(9, 5).Value = 14
$ws.Cells.Item(9, 6).Value = 0.276
$ws.Cells.Item(9, 7).Value = 0.426
$ws.Cells.Item(9, 9).Value = "MPM800"
$ws.Cells.Item(9, 10).Value = "Miscellaneous"
$ws.Cells.Item(9, 11).Value = "MPM800-1"
$ws.Cells.Item(9, 12).Value = 0.022
$ws.Cells.Item(9, 13).Value = 0.03
$ws.Cells.Item(9, 14).Value = 0.298
$ws.Cells.Item(9, 15).Value = 0.456
$ws.Cells.Item(9, 16).Value = "5V"
$ws.Cells.Item(9, 17).Value = 0.286
$ws.Cells.Item(9, 18).Value = 0.445

$ws.Range("A8:G8").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)
$ws.Range("I8:R8").Copy()
$ws.Range("I9:R9").PasteSpecial(-4122)

# C9 ("CPU Type") is left blank, but keeps a quote-prefixed left-aligned
# empty-text style (same style used by P8/L8/M8 -- style index 10).
$ws.Range("P8").Copy()
$ws.Cells.Item(9, 3).PasteSpecial(-4122)

$ws.Select()
$ws.Range("B8").Select()
